# Updated cryptos list on Wed Jul  5 19:28:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.529.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "'1.912.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'239.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("E8").Value = "  -2.57%  "

$ws.Range("D9").Value = "'0.06692"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("D10").Value = "'18.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "

$ws.Range("D11").Value = "'101.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.89%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07691"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.913.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("D14").Value = "'5.222"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "'0.6688"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "'30.541.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").Value = "'256.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.52%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").Value = "'0.000007474"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "

$ws.Range("D20").Value = "'12.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("D21").Value = "'5.385"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'6.292"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("D24").Value = "'9.329"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.64%  "

$ws.Range("D25").Value = "'166.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("E26").Value = "  -2.05%  "

$ws.Range("D27").Value = "'2.055"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.70%  "

$ws.Range("D28").Value = "'4.750"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "

$ws.Range("E29").Value = "  -2.60%  "

$ws.Range("D30").Value = "'1.382"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.67%  "

$ws.Range("D31").Value = "'1.514"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").Value = "'4.246"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("D33").Value = "'0.04714"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "

$ws.Range("D34").Value = "'0.7303"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.73%  "

$ws.Range("E35").Value = "  -3.74%  "

$ws.Range("D36").Value = "'0.9996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").Value = "'2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("D38").Value = "'0.01918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("D39").Value = "'2.617"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("D40").Value = "'6.233"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.33%  "

$ws.Range("D41").Value = "'74.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("D42").Value = "'1.968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.44%  "

$ws.Range("D43").Value = "'0.8624"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "

$ws.Range("D44").Value = "'105.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.29%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4235"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").Value = "'7.373"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.27%  "

$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("D49").Value = "'34.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "

$ws.Range("D50").Value = "'906.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.04%  "

$ws.Range("D51").Value = "'8.777"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "

